# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.057.21"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "'1.912.67"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.8330"
$ws.Range("E5").Value = "  +9.13%  "
$ws.Range("D6").Value = "'242.42"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.3252"
$ws.Range("E8").Value = "  +6.12%  "
$ws.Range("D9").Value = "'26.79"
$ws.Range("E9").Value = "  +4.93%  "
$ws.Range("D10").Value = "'0.07053"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("D11").Value = "'0.08041"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "'0.7539"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "'1.904.96"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "'5.240"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'93.01"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "'14.23"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "'30.048.43"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "'5.965"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "'245.53"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'0.000007786"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "'2.157.79"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'7.007"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "'0.1616"
$ws.Range("E25").Value = "  +24.95%  "
$ws.Range("D26").Value = "'169.77"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").Value = "'9.282"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").Value = "'18.99"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D29").Value = "'2.087"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").Value = "'1.372"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("D31").Value = "'1.520"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "'4.311"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "'0.05627"
$ws.Range("E33").Value = "  +6.51%  "
$ws.Range("D34").Value = "'4.107"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "'1.288"
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("D36").Value = "'0.7371"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'0.01919"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "'2.795"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "'0.4452"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").Value = "'6.019"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("D43").Value = "'0.8435"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D46").Value = "'7.627"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'101.27"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").Value = "'9.768"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'989.15"
$ws.Range("E49").Value = "  +9.45%  "
$ws.Range("D50").Value = "'2.063.55"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  +0.80%  "

# Rows 44/45 swapped order (RenderToken now ranks above PaxDollar) with new figures
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.905"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9993"
$ws.Range("E45").Value = "  -0.15%  "
